$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("1")
$ws.Range("F2").Value = 0.2321131094768888
$ws.Range("F3").Value = 0.2211114090295048
$ws.Range("F4").Value = 0.2545908439628514
$ws.Range("F5").Value = 0.292359944254258
$ws.Range("F7").Value = 0.2478714723361773
$ws.Range("F8").Value = 0.3421459761386313
$ws.Range("F9").Value = 0.3281994951562797
$ws.Range("F10").Value = 0.3233837799869357
$ws.Range("F11").Value = 0.3186070424084503
$ws.Range("F12").Value = 0.2915847476947995
$ws.Range("F13").Value = 0.3097520437677905
$ws.Range("F14").Value = 0.3045168564104406
$ws.Range("F15").Value = 0.2797147555867313
$ws.Range("F16").Value = 0.2815672720438374
$ws.Range("F17").Value = 0.3128949065232223
$ws.Range("F18").Value = 0.2843140643700365
$ws.Range("F19").Value = 0.3125600034729146

$ws = $wb.Worksheets.Item("3")
$ws.Range("F2").Value = 0.1811463416957392
$ws.Range("F3").Value = 0.1845024767365477
$ws.Range("F4").Value = 0.2976579084022617
$ws.Range("F5").Value = 0.3097550226805456
$ws.Range("F7").Value = 0.2543460058166472
$ws.Range("F8").Value = 0.3309530450481332
$ws.Range("F10").Value = 0.3155754231902388
$ws.Range("F11").Value = 0.2875124050411919
$ws.Range("F12").Value = 0.2706015633858573
$ws.Range("F13").Value = 0.3160795643620271
$ws.Range("F14").Value = 0.3126006577476584
$ws.Range("F15").Value = 0.3108005005884987
$ws.Range("F16").Value = 0.3147706543077511
$ws.Range("F17").Value = 0.2637268263214186
$ws.Range("F18").Value = 0.2846150939773607
$ws.Range("F19").Value = 0.2789936256639494

$ws = $wb.Worksheets.Item("5")
$ws.Range("F2").Value = 0.2194337381960671
$ws.Range("F3").Value = 0.1853323910528792
$ws.Range("F4").Value = 0.3382929116673225
$ws.Range("F7").Value = 0.291990374480973
$ws.Range("F8").Value = 0.4044551909186823
$ws.Range("F9").Value = 0.3312950345462412
$ws.Range("F10").Value = 0.3653705734394244
$ws.Range("F11").Value = 0.3156531724943581
$ws.Range("F12").Value = 0.2780800540824914
$ws.Range("F13").Value = 0.2857766086852866
$ws.Range("F14").Value = 0.2894708974557173
$ws.Range("F15").Value = 0.2863359343675875
$ws.Range("F16").Value = 0.3332879548132646
$ws.Range("F17").Value = 0.3178429090313465
$ws.Range("F18").Value = 0.2719375105021818
$ws.Range("F19").Value = 0.2916008267078748

$ws = $wb.Worksheets.Item("7")
$ws.Range("F2").Value = 0.1999344490641048
$ws.Range("F3").Value = 0.1680615406791622
$ws.Range("F4").Value = 0.298281665178119
$ws.Range("F7").Value = 0.2542916691963315
$ws.Range("F8").Value = 0.3374064063974575
$ws.Range("F9").Value = 0.3212737201289815
$ws.Range("F10").Value = 0.2955590949466989
$ws.Range("F11").Value = 0.2931946240593533
$ws.Range("F12").Value = 0.28781606854941
$ws.Range("F13").Value = 0.2733221574515677
$ws.Range("F14").Value = 0.2865468139027426
$ws.Range("F15").Value = 0.2426037527784093
$ws.Range("F16").Value = 0.2931400454222368
$ws.Range("F17").Value = 0.2878061442228339
$ws.Range("F19").Value = 0.3143198975649131

$ws = $wb.Worksheets.Item("Summary_All_Configs")
$ws.Range("C2").Value = 0.2904286895658676
$ws.Range("F2").Value = 0.03453674817974294
$ws.Range("I2").Value = 0.2211114090295048
$ws.Range("L2").Value = 0.3421459761386313
$ws.Range("C3").Value = 0.2821023196853642
$ws.Range("F3").Value = 0.04434134473086152
$ws.Range("I3").Value = 0.1811463416957392
$ws.Range("L3").Value = 0.3309530450481332
$ws.Range("C4").Value = 0.3017864309951214
$ws.Range("F4").Value = 0.05085246853242501
$ws.Range("I4").Value = 0.1853323910528792
$ws.Range("L4").Value = 0.4044551909186823
$ws.Range("C5").Value = 0.2790080264166018
$ws.Range("F5").Value = 0.04248480552246491
$ws.Range("I5").Value = 0.1680615406791622
$ws.Range("L5").Value = 0.3374064063974575

$ws = $wb.Worksheets.Item("Numeric_MEAN")
$ws.Range("C2").Value = 0.2904286895658676
$ws.Range("C3").Value = 0.2821023196853642
$ws.Range("C4").Value = 0.3017864309951214
$ws.Range("C5").Value = 0.2790080264166018

$ws = $wb.Worksheets.Item("Numeric_STD")
$ws.Range("C2").Value = 0.03453674817974294
$ws.Range("C3").Value = 0.04434134473086152
$ws.Range("C4").Value = 0.05085246853242501
$ws.Range("C5").Value = 0.04248480552246491

$ws = $wb.Worksheets.Item("Numeric_MIN")
$ws.Range("C2").Value = 0.2211114090295048
$ws.Range("C3").Value = 0.1811463416957392
$ws.Range("C4").Value = 0.1853323910528792
$ws.Range("C5").Value = 0.1680615406791622

$ws = $wb.Worksheets.Item("Numeric_MAX")
$ws.Range("C2").Value = 0.3421459761386313
$ws.Range("C3").Value = 0.3309530450481332
$ws.Range("C4").Value = 0.4044551909186823
$ws.Range("C5").Value = 0.3374064063974575
